$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0,0.1619081784284657,0.1945208102290912,-0.672860002190442,0.7034158221868738,1.570796280873402,-1.408888145739562),
    @(0.2073985369335523,0.1636540555887594,0.1944793538953593,-0.672975519928842,0.7033417592742593,1.570796281100031,-1.407142268608938),
    @(0.4147970738671046,0.1745421441903155,0.1942208132345191,-0.6736959413876683,0.7028798690720303,1.570796282513388,-1.396254180192416),
    @(0.6221956108006569,0.2002919558244658,0.1936093768695125,-0.675399703869848,0.7017875207142698,1.57079628585591,-1.370504368995864),
    @(0.8295941477342091,0.243380444329976,0.1925862288496451,-0.6782506975459399,0.6999596378417227,1.570796291449124,-1.327415881222608),
    @(1.036992684667762,0.3035824320851401,0.1911567160795255,-0.6822340233095296,0.6974057728506851,1.570796299263802,-1.267213894490529),
    @(1.244391221601314,0.3785110362998748,0.1893775157480047,-0.6871917506326238,0.6942271811418936,1.570796308990106,-1.192285291549146),
    @(1.451789758534866,0.4641580953078138,0.1873438027571155,-0.6928586754210456,0.6905938953694161,1.570796320107748,-1.106638233996711),
    @(1.659188295468418,0.5554345948584021,0.1851764171510118,-0.6988980778698286,0.6867217996895402,1.570796331956135,-1.015361735997295),
    @(1.866586832401971,0.6467110944089904,0.183009031544908,-0.7049374803186115,0.6828497040096643,1.570796343804521,-0.9240852379978779),
    @(2.073985369335523,0.7323581534169294,0.1809753185540188,-0.7106044051070333,0.6792164182371867,1.570796354922163,-0.8384381804454427),
    @(2.281383906269075,0.8072867576316641,0.179196118222498,-0.7155621324301276,0.6760378265283953,1.570796364648468,-0.7635095775040601),
    @(2.488782443202628,0.867488745386828,0.1777666054523784,-0.7195454581937172,0.6734839615373577,1.570796372463145,-0.703307590771981),
    @(2.69618098013618,0.9105772338923382,0.1767434574325111,-0.7223964518698091,0.6716560786648106,1.57079637805636,-0.6602191029987254),
    @(2.903579517069732,0.9363270455264889,0.1761320210675044,-0.7241002143519887,0.6705637303070501,1.570796381398882,-0.6344692918021734),
    @(3.110978054003284,0.9472151341280449,0.1758734804066642,-0.7248206358108151,0.670101840104821,1.570796382812239,-0.6235812033856507),
    @(3.318376590936837,0.9489610112883384,0.1758320240729323,-0.7249361535492151,0.6700277771922065,1.570796383038867,-0.6218353262550274)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($i + 1, $j + 1).Value = $row[$j]
    }
}
